# Localization sheet update: "Start" -> "Start Adventure" and a batch of
# new UI strings (death screen, settings, language picker, etc.) appended
# as rows 8-19. Column A is made to mirror column B (English) both in
# value and in font/style, matching the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Start Adventure"
$ws.Range("B2").Value = "Start Adventure"
$ws.Range("C2").Value = "开始冒险"
$ws.Range("A3").Value = "Options"
$ws.Range("B3").Value = "Options"
$ws.Range("C3").Value = "选项"
$ws.Range("A4").Value = "Quit"
$ws.Range("B4").Value = "Quit"
$ws.Range("C4").Value = "退出"
$ws.Range("A5").Value = "Volume"
$ws.Range("B5").Value = "Volume"
$ws.Range("C5").Value = "音量"
$ws.Range("A6").Value = "Back"
$ws.Range("B6").Value = "Back"
$ws.Range("C6").Value = "返回"
$ws.Range("A7").Value = "Home"
$ws.Range("B7").Value = "Home"
$ws.Range("C7").Value = "主页"
$ws.Range("A8").Value = "You Died"
$ws.Range("B8").Value = "You Died"
$ws.Range("C8").Value = "阵亡"
$ws.Range("A9").Value = "Play Again"
$ws.Range("B9").Value = "Play Again"
$ws.Range("C9").Value = "再次游玩"
$ws.Range("A10").Value = "Return Home"
$ws.Range("B10").Value = "Return Home"
$ws.Range("C10").Value = "返回主页"
$ws.Range("A11").Value = "Graphics: "
$ws.Range("B11").Value = "Graphics: "
$ws.Range("C11").Value = "图像"
$ws.Range("A12").Value = "Volume:"
$ws.Range("B12").Value = "Volume:"
$ws.Range("C12").Value = "音量"
$ws.Range("A13").Value = "Language:"
$ws.Range("B13").Value = "Language:"
$ws.Range("C13").Value = "语言"
$ws.Range("A14").Value = "Settings"
$ws.Range("B14").Value = "Settings"
$ws.Range("C14").Value = "设置"
$ws.Range("A15").Value = "Deep Dive Descent"
$ws.Range("B15").Value = "Deep Dive Descent"
$ws.Range("C15").Value = "深浅迷航"
$ws.Range("A16").Value = "Language"
$ws.Range("B16").Value = "Language"
$ws.Range("C16").Value = "语言"
$ws.Range("A17").Value = "Return"
$ws.Range("B17").Value = "Return"
$ws.Range("C17").Value = "返回"
$ws.Range("A18").Value = "English"
$ws.Range("B18").Value = "English"
$ws.Range("C18").Value = "English"
$ws.Range("A19").Value = "Chinese"
$ws.Range("B19").Value = "Chinese"
$ws.Range("C19").Value = "中文"

# Column A had no explicit style before; give it the same "Noto Sans"
# font style already used by column B (style index 1) by copying format
# instead of touching Font.Name directly (which would mint a duplicate font).
$ws.Range("B2:B19").Copy()
$ws.Range("A2:A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the selection to the last-edited cell, like the real edit session did.
$ws.Range("C19").Select() | Out-Null
